$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.837.31'
$ws.Range('E2').Value = '  +1.65%  '
$ws.Range('D3').Value = '3.633.97'
$ws.Range('E3').Value = '  +3.79%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '604.76'
$ws.Range('E5').Value = '  +0.16%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '199.96'
$ws.Range('E6').Value = '  +2.59%  '
$ws.Range('E7').Value = '  +0.26%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('E9').Value = '  +9.47%  '
$ws.Range('E10').Value = '  -0.62%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '53.84'
$ws.Range('E11').Value = '  +0.88%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000306'
$ws.Range('E12').Value = '  +2.07%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '9.56'
$ws.Range('E13').Value = '  +0.95%  '
$ws.Range('D14').Value = '4.210.05'
$ws.Range('E14').Value = '  +3.62%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '641.05'
$ws.Range('E15').Value = '  +7.96%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '13.00'
$ws.Range('E16').Value = '  +1.80%  '
$ws.Range('D17').Value = '70.937.99'
$ws.Range('E17').Value = '  +1.58%  '
$ws.Range('D18').Value = '3.622.13'
$ws.Range('E18').Value = '  +3.31%  '
$ws.Range('E19').Value = '  +0.39%  '
$ws.Range('E20').Value = '  +0.71%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.999'
$ws.Range('E21').Value = '  +1.27%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '18.34'
$ws.Range('E22').Value = '  +1.57%  '
$ws.Range('E23').Value = '  +1.53%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '104.16'
$ws.Range('E24').Value = '  +1.92%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.63'
$ws.Range('E25').Value = '  -0.57%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.00'
$ws.Range('E26').Value = '  -4.84%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.49'
$ws.Range('E27').Value = '  -3.30%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.74'
$ws.Range('E28').Value = '  +2.14%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '33.96'
$ws.Range('E29').Value = '  +2.10%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.78'
$ws.Range('E30').Value = '  +12.61%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.19'
$ws.Range('E31').Value = '  +2.20%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '12.25'
$ws.Range('E32').Value = '  -1.36%  '
$ws.Range('E33').Value = '  +1.06%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '63.35'
$ws.Range('E34').Value = '  +0.44%  '
$ws.Range('B35').Value = 'Maker'
$ws.Range('C35').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D35').Value = '4.032.61'
$ws.Range('E35').Value = '  +8.64%  '
$ws.Range('B36').Value = 'PEPE'
$ws.Range('C36').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D36').Value = '0.0₃0879'
$ws.Range('E36').Value = '  +6.13%  '
$ws.Range('E37').Value = '  +0.37%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '517.48'
$ws.Range('E38').Value = '  +9.98%  '
$ws.Range('E39').Value = '  -2.06%  '
$ws.Range('B40').Value = 'InjectiveProtocol'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '36.72'
$ws.Range('E40').Value = '  +0.80%  '
$ws.Range('B41').Value = 'TheGraph'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.390'
$ws.Range('E41').Value = '  -0.48%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.54'
$ws.Range('E42').Value = '  -2.82%  '
$ws.Range('E43').Value = '  +2.37%  '
$ws.Range('E44').Value = '  +2.23%  '
$ws.Range('E45').Value = '  +6.80%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.94'
$ws.Range('E46').Value = '  +4.56%  '
$ws.Range('E47').Value = '  +0.70%  '
$ws.Range('E48').Value = '  +2.60%  '
$ws.Range('E49').Value = '  -0.21%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.000250'
$ws.Range('E50').Value = '  +2.37%  '
$ws.Range('E51').Value = '  +1.39%  '
